# Generate Report for Handback
# Update the handoff/handback timestamps recorded for file
# "38977caf-9a6f-41cd-b1c2-99c801909518" (the report row that is already
# in sync with en-US) across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 38977caf... row (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-05 16:53:07"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) for the 38977caf... row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-05 16:52:58"
$wsZhCn.Range("K3").Value = "2016-09-05 16:53:31"

# de-de sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) for the 38977caf... row (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-05 16:53:07"
$wsDeDe.Range("K3").Value = "2016-09-05 16:53:39"
